$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the scraped cryptocurrency Price / Volume(1h) figures (a handful of rows
# also swap Coin/Link because the source ranking order shifted), matching the
# GitHub Actions commit "Updated cryptos list on Thu Mar 30 08:28:44 UTC 2023".

# Price column ("D") holds plain text even when it looks numeric (e.g. "1.001"),
# so force Text format first -- otherwise Excel would silently coerce values like
# "316.61" into real numbers and drop the original formatting/trailing zeros.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.721.67"
$ws.Range("E2").Value = "  +1.29%  "

$ws.Range("D3").Value = "1.807.88"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.61%  "

$ws.Range("D5").Value = "316.61"
$ws.Range("E5").Value = "  -0.26%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.51%  "

$ws.Range("D7").Value = "0.5360"
$ws.Range("E7").Value = "  -5.25%  "

$ws.Range("D8").Value = "0.3792"
$ws.Range("E8").Value = "  -1.47%  "

$ws.Range("D9").Value = "0.07536"
$ws.Range("E9").Value = "  -1.46%  "

$ws.Range("D10").Value = "42.68"
$ws.Range("E10").Value = "  -1.18%  "

$ws.Range("D11").Value = "1.119"
$ws.Range("E11").Value = "  -1.66%  "

$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.72%  "

$ws.Range("D13").Value = "20.95"
$ws.Range("E13").Value = "  -2.00%  "

$ws.Range("D14").Value = "6.194"
$ws.Range("E14").Value = "  -0.82%  "

$ws.Range("D15").Value = "7.391"
$ws.Range("E15").Value = "  +2.40%  "

$ws.Range("D16").Value = "1.805.43"
$ws.Range("E16").Value = "  +0.34%  "

$ws.Range("D17").Value = "90.65"
$ws.Range("E17").Value = "  -1.78%  "

$ws.Range("E18").Value = "  -1.42%  "

$ws.Range("D19").Value = "0.06460"
$ws.Range("E19").Value = "  -0.80%  "

$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.48%  "

$ws.Range("D21").Value = "17.29"
$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("D22").Value = "5.918"
$ws.Range("E22").Value = "  -1.20%  "

$ws.Range("D23").Value = "28.737.62"
$ws.Range("E23").Value = "  +1.44%  "

$ws.Range("D24").Value = "11.21"
$ws.Range("E24").Value = "  -0.69%  "

$ws.Range("D25").Value = "2.107"
$ws.Range("E25").Value = "  +0.42%  "

$ws.Range("D26").Value = "160.99"
$ws.Range("E26").Value = "  +2.84%  "

$ws.Range("D27").Value = "20.50"
$ws.Range("E27").Value = "  -1.50%  "

$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D28").Value = "2.016.27"
$ws.Range("E28").Value = "  +0.34%  "

$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "2.365"
$ws.Range("E29").Value = "  -0.65%  "

$ws.Range("D30").Value = "123.18"
$ws.Range("E30").Value = "  -0.30%  "

$ws.Range("D31").Value = "1.108"
$ws.Range("E31").Value = "  -3.98%  "

$ws.Range("D32").Value = "0.1061"
$ws.Range("E32").Value = "  +1.65%  "

$ws.Range("D33").Value = "5.668"
$ws.Range("E33").Value = "  -1.64%  "

$ws.Range("E34").Value = "  +2.13%  "

$ws.Range("D35").Value = "0.2263"
$ws.Range("E35").Value = "  +6.05%  "

$ws.Range("D36").Value = "0.06448"
$ws.Range("E36").Value = "  +6.15%  "

$ws.Range("D37").Value = "0.02316"
$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("D38").Value = "8.776"
$ws.Range("E38").Value = "  +1.01%  "

$ws.Range("D39").Value = "5.052"
$ws.Range("E39").Value = "  +0.27%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "1.215"
$ws.Range("E40").Value = "  +5.35%  "

$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "11.31"
$ws.Range("E41").Value = "  -3.20%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.6257"
$ws.Range("E42").Value = "  -2.73%  "

$ws.Range("D43").Value = "0.9997"
$ws.Range("E43").Value = "  +0.39%  "

$ws.Range("D44").Value = "1.399"
$ws.Range("E44").Value = "  +0.41%  "

$ws.Range("E45").Value = "  -0.75%  "

$ws.Range("D46").Value = "0.5883"
$ws.Range("E46").Value = "  -1.91%  "

$ws.Range("D47").Value = "3.689"
$ws.Range("E47").Value = "  -0.04%  "

$ws.Range("D48").Value = "126.21"
$ws.Range("E48").Value = "  +3.24%  "

$ws.Range("D49").Value = "1.958"
$ws.Range("E49").Value = "  +1.00%  "

$ws.Range("D50").Value = "1.153"
$ws.Range("E50").Value = "  +0.94%  "

$ws.Range("D51").Value = "0.06889"
$ws.Range("E51").Value = "  +0.83%  "

# Drop the explicit Text number-format override again so the cells keep matching
# the original (unstyled) price cells once the values are safely stored as text.
$priceRange.Style = "Normal"